$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting existing columns B..AA to C..AB
$ws.Columns("B:B").Insert(4)
$ws.Columns("B:B").ColumnWidth = 20.14

# New snapshot timestamp goes into the header row
$ws.Range("B1").Value = "2025-12-22 22:24"

# Carry forward the most recent known price into the new column for each
# product row (same as what used to be in the old column B, now shifted to C),
# except for rows where the product had no price recorded this time around.
$ws.Range("B2").Value = 929
$ws.Range("B3").Value = 569
$ws.Range("B4").Value = 299
$ws.Range("B5").Value = 569
$ws.Range("B6").Value = 499
$ws.Range("B7").Value = 569
$ws.Range("B8").Value = 929
$ws.Range("B9").Value = 299
$ws.Range("B11").Value = 2997
$ws.Range("B12").Value = 569
$ws.Range("B13").Value = 569
$ws.Range("B15").Value = 499
$ws.Range("B16").Value = 299
$ws.Range("B17").Value = 929
$ws.Range("B18").Value = 499
$ws.Range("B19").Value = 1299
$ws.Range("B20").Value = 929
$ws.Range("B21").Value = 499
$ws.Range("B22").Value = 299
$ws.Range("B23").Value = 1299
$ws.Range("B24").Value = 929
$ws.Range("B25").Value = 929
$ws.Range("B26").Value = 1299
# B10 and B14 are left blank (no price recorded for this snapshot)

Write-Host "done"
